# Applies the "16.11.2018" work-log entry to the "Eetu Pihamäki" worksheet
# row 31, matching the commit that added this row's data + a new shared
# string describing the work done that day.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Eetu Pihamäki")

# Fill in the new time-tracking entry for 16.11.2018 (row 31):
#   A31 = date, B31 = start time, C31 = end time (stored as the underlying
#   serial/fraction numbers, same as the other rows in this table),
#   D31 already holds a shared formula (C31-B31) that recalculates,
#   E31 = sprint number, F31 = task description text.
$ws.Range("A31").Value2 = 43420
$ws.Range("B31").Value2 = 0.45833333333333331
$ws.Range("C31").Value2 = 0.5
$ws.Range("E31").Value2 = 4
$ws.Range("F31").Value = "1h Samaa kuin viime kerralla. https://github.com/Eetu95/Open-source-IdM-solution/blob/master/Eetun%20muistiinpanoja/Ty%C3%B6t%20-%2016.11.2018.txt"

# Row 31 now wraps a full task description like the other populated rows, so
# its height grows to match them (same 60pt row height used by similar rows).
$ws.Rows.Item(31).RowHeight = 60

# Update the active selection to match the new edit location.
$ws.Range("F31").Select()

$wb.Save()
